# Update Stock (column D) quantities on the "Inventory" sheet and
# move the active selection, per the "seperated shopping sequence from main" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 53
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 64
$ws.Range("D7").Value = 60
$ws.Range("D10").Value = 1

$ws.Activate()
$ws.Range("M8").Select()
